$d = $word.ActiveDocument

# --- Edit 1: first paragraph -------------------------------------------------
# Original: "This is a Microsoft word document."
# New: "This is a Microsoft word document.  " (plain run, two trailing spaces)
#      + red run "(This is a change – Ve"
#      + red run "rsion for main branch"
#      + red run ")"

$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "This is a Microsoft word document.  ", 2)

$p1 = $d.Paragraphs(1).Range
$insPos = $p1.End - 1

$r1 = $d.Range($insPos, $insPos)
$r1.InsertAfter([char]40 + "This is a change " + [char]0x2013 + " Ve")
$r1.Font.Color = 255

$insPos = $insPos + $r1.End - $r1.Start - ($r1.End - $r1.Start) + ($r1.End - $r1.Start)
$insPos = $r1.End
$r2 = $d.Range($insPos, $insPos)
$r2.InsertAfter("rsion for main branch")
$r2.Font.Color = 255

$insPos = $r2.End
$r3 = $d.Range($insPos, $insPos)
$r3.InsertAfter(")")
$r3.Font.Color = 255

# --- Edit 2: remove the last paragraph ("ank God almighty, we are free at last.") ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastPara.Range.Delete()
